$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 354 (shifting existing rows 354..443 down to 355..444).
$ws.Rows.Item(354).Insert()

# Populate the newly inserted row with the new record.
$ws.Range("A354").Value = 8
$ws.Range("B354").Value = "Terminal La Palmera de La Serena"
$ws.Range("C354").Value = "Coquimbo"
$ws.Range("D354").Value = 44754
$ws.Range("E354").Value = 4
$ws.Range("F354").Value = 100114001
$ws.Range("G354").Value = "Papa"
$ws.Range("H354").Value = "Asterix"
$ws.Range("I354").Value = "1a (guarda)"
$ws.Range("J354").Value = 2460
$ws.Range("K354").Value = 9500
$ws.Range("L354").Value = 10000
$ws.Range("M354").Value = 9750
$ws.Range("N354").Value = "$/saco 25 kilos"
$ws.Range("O354").Value = "Región de Los Lagos"
$ws.Range("P354").Value = 390
$ws.Range("Q354").Value = 25
$ws.Range("R354").Value = "Hortaliza"
